$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets("Config")

# 1) B3 dropdown: "Black List" -> "All Students" (recalculates C3 automatically)
$ws.Range("B3").Value = "All Students"

# 2) Append a new row 12 (Title / 2nd MidTerm), copying the formatting of row 11
#    so the new cells inherit the same style as the rest of the key/value list.
$ws.Range("A11").Copy()
$ws.Range("A12:B12").PasteSpecial(-4122)
$ws.Range("A12").Value = "Title"
$ws.Range("B12").Value = "2nd MidTerm"

# 3) Update the view: scroll to show the new rows, move the selection to B10
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("B10").Select()

# 4) Remove the blank extra worksheet "Sheet1"
$wb.Sheets("Sheet1").Delete()
